# Append 21 new daily NAV rows (2024-08-28 .. 2024-09-25) to the end of
# the "Luxury_Goods" data sheet, continuing directly after the existing
# row 610 (2024-09-25) so the sheet grows from A1:J610 to A1:J631.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Date, BLUESTARCO.NS, METROBRAND.NS, VOLTAS.NS, WHIRLPOOL.NS,
# CROMPTON.NS, Basket Value, daily return, NAV
$newRows = @(
    ,@("2024-08-28", 1746.650024414062, 1322.349975585938, 1771.75, 2073.949951171875, 462.7999877929688, 8303.099914550781, 0, 217.5326560482892)
    ,@("2024-08-29", 1695.900024414062, 1307.849975585938, 1788.400024414062, 2194.800048828125, 465, 8381.950073242188, 0.009496472342001468, 219.5984488999339)
    ,@("2024-08-30", 1703.900024414062, 1306.050048828125, 1743.650024414062, 2183.5, 477.0499877929688, 8368.250061035156, -0.001634465975974491, 219.2395227068302)
    ,@("2024-09-02", 1660.599975585938, 1280.449951171875, 1769.650024414062, 2220.5, 464.8500061035156, 8325.749969482422, -0.005078731065963999, 218.1260641319719)
    ,@("2024-09-03", 1683.75, 1248.449951171875, 1810.949951171875, 2205.5, 469.6000061035156, 8357.449920654297, 0.003807458942205738, 218.9565701653793)
    ,@("2024-09-04", 1699.5, 1288.900024414062, 1780.25, 2217.949951171875, 466.5499877929688, 8386.249938964844, 0.003446029420932761, 219.7111009480757)
    ,@("2024-09-05", 1694.699951171875, 1272.300048828125, 1783.150024414062, 2236.14990234375, 467.75, 8389.549926757812, 0.0003934998142180441, 219.7975572254804)
    ,@("2024-09-06", 1724.449951171875, 1234.300048828125, 1778.650024414062, 2224.39990234375, 463.3999938964844, 8351.999908447266, -0.004475808432915339, 218.8137854653164)
    ,@("2024-09-09", 1789.300048828125, 1234.300048828125, 1817.949951171875, 2138.60009765625, 460.7999877929688, 8362.550109863281, 0.001263194627833399, 219.0901898636121)
    ,@("2024-09-10", 1799.949951171875, 1249.650024414062, 1830.099975585938, 2140.14990234375, 466.9500122070312, 8420.699890136719, 0.006953594239734628, 220.6136541458301)
    ,@("2024-09-11", 1788.25, 1255, 1828.699951171875, 2119.85009765625, 459.2999877929688, 8369.700012207031, -0.006056489198650145, 219.2775099324211)
    ,@("2024-09-12", 1816.650024414062, 1257.449951171875, 1852.949951171875, 2111, 459, 8415.049926757812, 0.005418344084571652, 220.4656309312431)
    ,@("2024-09-13", 1888, 1264.349975585938, 1921.550048828125, 2101.35009765625, 451.75, 8530.500122070312, 0.01371949023681921, 223.4903070023585)
    ,@("2024-09-16", 1934.900024414062, 1259.75, 1916, 2047.199951171875, 447.8999938964844, 8501.54995727539, -0.003393724210849177, 222.7318425365945)
    ,@("2024-09-17", 1891.199951171875, 1257.550048828125, 1904.050048828125, 2034.699951171875, 452.75, 8445.75, -0.006563504014657771, 221.2699411939134)
    ,@("2024-09-18", 1879.449951171875, 1232.050048828125, 1889.400024414062, 2042.550048828125, 448.8500061035156, 8390.000091552734, -0.006600942302017657, 219.8093510789215)
    ,@("2024-09-19", 1866.650024414062, 1234.5, 1929.199951171875, 2025.699951171875, 450, 8406.049926757812, 0.00191297199403341, 220.2298402115621)
    ,@("2024-09-20", 1930.099975585938, 1246.550048828125, 1928.400024414062, 2062.300048828125, 446.7999877929688, 8507.750061035156, 0.01209844518691423, 222.8942788618846)
    ,@("2024-09-23", 2012.849975585938, 1229.900024414062, 1882.449951171875, 2039, 445.75, 8501.449951171875, -0.000740514215636784, 222.7292224798032)
    ,@("2024-09-24", 2026, 1238.800048828125, 1914.400024414062, 2064.550048828125, 439.25, 8561.500122070312, 0.007063521075032611, 224.302475036815)
    ,@("2024-09-25", 2051.39990234375, 1252.949951171875, 1910.849975585938, 2112.050048828125, 441.5499877929688, 8651.899841308594, 0.01055886444540762, 226.6708544654981)
)

$startRow = 611
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds the date as literal text (matches the existing
    # "yyyy-MM-dd" text entries above it) rather than a real Excel date
    # serial, so force text entry with a leading apostrophe, then strip
    # the formatting stamp that the quote-prefix leaves behind.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 1).ClearFormats()

    $ws.Cells.Item($r, 3).Value  = $row[1]   # BLUESTARCO.NS
    $ws.Cells.Item($r, 4).Value  = $row[2]   # METROBRAND.NS
    $ws.Cells.Item($r, 5).Value  = $row[3]   # VOLTAS.NS
    $ws.Cells.Item($r, 6).Value  = $row[4]   # WHIRLPOOL.NS
    $ws.Cells.Item($r, 7).Value  = $row[5]   # CROMPTON.NS
    $ws.Cells.Item($r, 8).Value  = $row[6]   # Basket Value
    $ws.Cells.Item($r, 9).Value  = $row[7]   # daily return
    $ws.Cells.Item($r, 10).Value = $row[8]   # NAV
}
